# Update investment cost results pulled from server for sheets 2025, 2030, 2035

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 1320.21302159999
$ws.Range("E2").Value = 294386.0104869407
$ws.Range("G2").Value = 80959.25712661834
$ws.Range("I2").Value = 141747.6361456303
$ws.Range("L2").Value = 525583.2829870571
$ws.Range("M2").Value = 111324.924271
$ws.Range("N2").Value = 71365.52612246884
$ws.Range("O2").Value = 68011.00857379404

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 46114.99750964541
$ws.Range("E2").Value = 275102.3308773747
$ws.Range("I2").Value = 186257.2480463011
$ws.Range("L2").Value = 296535.5563106379
$ws.Range("M2").Value = 104257.9211117674
$ws.Range("N2").Value = 37190.87036815702
$ws.Range("O2").Value = 27504.58861597344

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 19901.87149144481
$ws.Range("B2").Value = 20277.90980329371
$ws.Range("E2").Value = 116457.919756194
$ws.Range("I2").Value = 166252.3492318689
$ws.Range("M2").Value = 53692.05352104004
$ws.Range("N2").Value = 50073.96144272469
$ws.Range("O2").Value = 56800.26050142136
